$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The edit: insert a new row (row 6) above the data table containing the
# note "Sources found via https://www.data.va.gov". Everything that was at
# row 7 (table header) and below shifts down by one row.
#
# Because the automated row-insert does not carry the worksheet's existing
# hyperlink objects down with it, we first remove them and recreate them
# afterwards, pointed at their (now shifted) cells with their original
# target URLs, so the final layout/relationships match.
# ---------------------------------------------------------------------------

# Remember the existing hyperlink target URLs, in the order they appear
# (D8,E8,D9,E9,...,D13,E13) before we touch anything.
$hyperlinkTargets = @(
  "https://www.data.va.gov/stories/s/4tx7-hu2d",
  "https://www.data.va.gov/api/views/dwpj-hgp7/rows.csv?date=20230707&accessType=DOWNLOAD",
  "https://www.data.va.gov/dataset/VA-Opioid-Prescribing-Facilities/at5r-w2x9",
  "https://www.data.va.gov/api/views/dwpj-hgp7/rows.csv?accessType=DOWNLOAD",
  "https://www.data.va.gov/dataset/Opioid-Prescribing-Rates-at-VA-Facilities-2012-201/dwpj-hgp7",
  "https://www.data.va.gov/api/views/dwpj-hgp7/rows.csv?accessType=DOWNLOAD&bom=true&format=true",
  "https://www.data.va.gov/dataset/Department-of-Veterans-Affairs-Opioid-Dispensing-D/qvgv-ry3b",
  "https://www.data.va.gov/download/qvgv-ry3b/text%2Fplain",
  "https://www.data.va.gov/dataset/Department-of-Veterans-Affairs-Opioid-Dispensing-D/9478-kz49",
  "https://www.data.va.gov/download/9478-kz49/text%2Fplain",
  "https://www.data.va.gov/dataset/Department-of-Veterans-Affairs-Opioid-Dispensing-D/um24-98en",
  "https://www.data.va.gov/download/um24-98en/text%2Fplain"
)

# The cells that will hold those same hyperlinks once the new row has been
# inserted (each original row number + 1).
$hyperlinkCells = @("D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14")

# Drop the current hyperlinks so stale cell references are not left behind
# once the rows shift.
$ws.Hyperlinks.Delete()

# Insert the new row above the current row 6 (the blank separator row),
# pushing the header/table down by one row.
$ws.Range("A6:E6").EntireRow.Insert()

# Fill in the new note row.
$ws.Range("A6").Value = "Sources found via https://www.data.va.gov"

# The three "Department ... Opioid ... Data" rows directly under the header
# had a leftover, visually inert style on column C; that style is gone after
# the edit, so clear the (no-op) formatting back to the default.
$ws.Range("C9:C11").ClearFormats()

# Recreate the hyperlinks against their shifted cells with their original
# target URLs.
for ($i = 0; $i -lt $hyperlinkCells.Length; $i++) {
  [void]$ws.Hyperlinks.Add($ws.Range($hyperlinkCells[$i]), $hyperlinkTargets[$i])
}

# Recreating the hyperlinks can introduce a stray duplicate style; make sure
# the affected cells keep using the normal "Hyperlink" cell style.
$ws.Range("D9:E14").Style = "Hyperlink"

# Match the saved selection (cell A7 was active when the author last saved).
[void]$ws.Range("A7").Select()
